$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Sewer Tentacles" horror card as row 10 (non-description columns first)
$ws.Range("A10").Value = "horror"
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "Sewer Tentacles"
$ws.Range("D10").Value = 4
$ws.Range("E10").Value = ":H: :H:"

# Update the description text on row 9 (F9): "discovered" -> "revealed"
$ws.Range("F9").Value = "If revealed by a player that's not in the same space, flip this face down."

# Description for the new row
$ws.Range("F10").Value = "When revealed, all players on the edge of the city lose 2 :heart:"

# Match formatting (wrap text) used by the rest of the table, and give the
# new card row the correct height so it renders as a proper square card.
$ws.Range("A10:F10").WrapText = $true
$ws.Rows.Item(10).RowHeight = 45

# Update selection to reflect where the user ended up editing
$ws.Range("F10").Select()
